$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8027499914169312
$ws.Range("C2").Value = 0.6780000925064087
$ws.Range("D2").Value = 0.4382500350475311

# Row 3
$ws.Range("B3").Value = 0.8540000319480896
$ws.Range("C3").Value = 0.8087500333786011
$ws.Range("D3").Value = 0.7212499976158142

# Row 4
$ws.Range("B4").Value = 0.7913309335708618
$ws.Range("C4").Value = 0.8316702246665955
$ws.Range("D4").Value = 0.3585298657417297

# Row 5
$ws.Range("B5").Value = 0.8216931223869324
$ws.Range("C5").Value = 0.8605567216873169
$ws.Range("D5").Value = 0.6356503963470459

# Row 6 (D6 newly added)
$ws.Range("B6").Value = 0.8035541772842407
$ws.Range("C6").Value = 0.675298810005188
$ws.Range("D6").Value = 0.4656227827072144

# Row 7 (D7 newly added)
$ws.Range("B7").Value = 0.8587315082550049
$ws.Range("C7").Value = 0.8079832792282104
$ws.Range("D7").Value = 0.7507182359695435

# Row 8
$ws.Range("B8").Value = 0.8720000982284546
$ws.Range("C8").Value = 0.9665000438690186
$ws.Range("D8").Value = 0.8567500114440918

# Row 9
$ws.Range("B9").Value = 0.874750018119812
$ws.Range("C9").Value = 0.9662500619888306
$ws.Range("D9").Value = 0.7445000410079956

# Row 10 (D10 newly added)
$ws.Range("B10").Value = 0.8234329223632812
$ws.Range("C10").Value = 0.5971059799194336
$ws.Range("D10").Value = 0.4986208081245422

# Row 11 (B11 updated, C11 and D11 cleared/removed)
$ws.Range("B11").Value = 0.8581515550613403
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()

# Row 12 (D12 newly added)
$ws.Range("B12").Value = 0.8234329223632812
$ws.Range("C12").Value = 0.5971059799194336
$ws.Range("D12").Value = 0.4986208081245422
